$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G9").Value = 1.91
$ws.Range("I9").Value = 4.5
$ws.Range("Q9").Value = 2.2
$ws.Range("R9").Value = 1.65
$ws.Range("U9").Value = 2
$ws.Range("V9").Value = 1.75
$ws.Range("W9").Value = 6
$ws.Range("X9").Value = 8
$ws.Range("Z9").Value = 15
$ws.Range("AC9").Value = 7.5
$ws.Range("AG9").Value = 401
$ws.Range("AH9").Value = 11
$ws.Range("AK9").Value = 51
$ws.Range("AU9").Value = 9
$ws.Range("BB9").Value = 126
